# docs: align Project.xlsx with current workflow
#
# Rebuilds the "Initial" sheet's header/guidance block to match the
# current AI-assisted hearing -> task-list / man-hour workflow:
#   - "Hearing" / "AI Output" banner row (merged, bold, centered)
#   - Input / Output / Remarks / Task List (AI) / Man Hour Estimation (AI)
#     column headers (bold, centered, wrapped)
#   - Updated Indonesian guidance text in row 3, now including a
#     Remarks column note and a Task List formula hint
#   - A new row 4 note describing how the AI output spills into D:E
#   - The stale placeholder rows that used to live in A4/B4 are cleared

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Initial"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Initial")
$ws.Visible = -1

# --- Merge the banner cells first, so later formatting only has to touch
#     the resulting single anchor cells (A1 / D1) instead of spilling
#     style onto every cell the merge swallows (B1, C1, E1). ---------------
$ws.Range("A1:C1").Merge()
$ws.Range("D1:E1").Merge()

# --- Row 1: banner -----------------------------------------------------
$ws.Range("A1").Value = "Hearing"
$ws.Range("D1").Value = "AI Output"

# --- Row 2: column headers ----------------------------------------------
$ws.Range("A2").Value = "Input"
$ws.Range("B2").Value = "Output"
$ws.Range("C2").Value = "Remarks"
$ws.Range("D2").Value = "Task List (AI)"
$ws.Range("E2").Value = "Man Hour Estimation (AI)"

# --- Row 3: guidance text (updated wording, now also covering Remarks
#     and the Task List formula) -----------------------------------------
$ws.Range("A3").Value = "Isi ringkasan kebutuhan / input (dari hearing)."
$ws.Range("B3").Value = "Isi deliverables / output yang diharapkan (dari hearing)."
$ws.Range("C3").Value = "Isi catatan: constraint, asumsi, dan pertanyaan terbuka."
$ws.Range("D3").Value = "Di Google Sheets, masukkan formula di kolom Task List: =GENERATE_ESTIMATION(Ax,Bx,Cx)"

# --- Row 4: old placeholder text is gone; replaced by a single spill note
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("D4").Value = "Output akan mengisi (spill) kolom D:E sebagai baris [Task, Hours]."

# --- Formatting: banner + header cells (bold, centered, wrapped) ---------
# Only the cells that actually carry text get the style. A multi-area
# union range here only reliably applies formatting to its first area, so
# walk `.Areas` and format each one individually -- this leaves blank
# neighbours (already absorbed into the A1:C1 / D1:E1 merges) untouched.
$headerCells = $ws.Range("A1,D1,A2,B2,C2,D2,E2")
foreach ($area in $headerCells.Areas) {
    $area.WrapText = $true
    $area.VerticalAlignment = -4108     # xlCenter
    $area.HorizontalAlignment = -4108   # xlCenter
    $area.Font.Bold = $true
}

# --- Formatting: guidance cells (left, top, wrapped) ------------------------
$bodyCells = $ws.Range("A3,B3,C3,D3,D4")
foreach ($area in $bodyCells.Areas) {
    $area.WrapText = $true
    $area.VerticalAlignment = -4160     # xlTop
    $area.HorizontalAlignment = -4131   # xlLeft
}

$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "FS" -- content unchanged, just make sure it stays visible
# ---------------------------------------------------------------------
$fs = $wb.Worksheets.Item("FS")
$fs.Visible = -1
